# Sprint4 backlog update
# - Row 9 (B9) user story renamed from "Profile view and Edit page" to "Sign in page"
# - The two now-empty trailing rows (12 & 13) are removed so the sheet's
#   used range shrinks back down to A1:F11
# - Selection is left on B10, matching the saved workbook's cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Sign in page"

$ws.Rows.Item(12).Delete() | Out-Null
$ws.Rows.Item(12).Delete() | Out-Null

$ws.Range("B10").Select() | Out-Null
